$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$result = $find.Execute("ach.label and ach.text %}: {% endif %}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ach.label and ach.text %} {% endif %}", 2)

Write-Output "Replace result: $result"
